$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample name shared string (used by A4:A9) from "248F-19 " to "SOX 717-C 1 "
$ws.Range("A4:A9").Value = "SOX 717-C 1 "

# Update existing rows 4-6 with new data
$ws.Range("B4").Value = 20201209059
$ws.Range("C4").Value = 159.6991324768351
$ws.Range("D4").Value = 0.1049552903810511
$ws.Range("E4").Value = 46.20157618021523
$ws.Range("F4").Value = 1.462899470367886

$ws.Range("B5").Value = 20201209060
$ws.Range("C5").Value = 159.6154088578845
$ws.Range("D5").Value = 0.08121870658997199
$ws.Range("E5").Value = 37.48976383711215
$ws.Range("F5").Value = 0.8996805363756503

$ws.Range("B6").Value = 20201209061
$ws.Range("C6").Value = 158.960571292258
$ws.Range("D6").Value = 0.1339382949238062
$ws.Range("E6").Value = 42.1683928693664
$ws.Range("F6").Value = 1.489555878658869

# Add new rows 7-9 (Sample Name already set above via A4:A9)
$ws.Range("B7").Value = 20201209062
$ws.Range("C7").Value = 160.0310128606721
$ws.Range("D7").Value = 0.1143225937587857
$ws.Range("E7").Value = 32.34361313266743
$ws.Range("F7").Value = 0.802813239473873

$ws.Range("B8").Value = 20201209063
$ws.Range("C8").Value = 160.2176175865996
$ws.Range("D8").Value = 0.08067987777175896
$ws.Range("E8").Value = 36.82344379992367
$ws.Range("F8").Value = 0.8276634857217821

$ws.Range("B9").Value = 20201209064
$ws.Range("C9").Value = 159.9232372688059
$ws.Range("D9").Value = 0.0793632401463726
$ws.Range("E9").Value = 37.35252521284222
$ws.Range("F9").Value = 0.8666815264503328
